# Weekly update: insert two new price rows for "Terminal La Palmera de La
# Serena - Cebolla" (row 828/829), pushing the existing historical rows
# down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 828:829 - this shifts the former row 828 onward
# down by two rows (old 828 -> new 830, ..., old 881 -> new 883) and keeps
# the date-formatted style of column D for the new rows.
$ws.Rows("828:829").Insert()

# New row 828 data
$ws.Cells.Item(828,1).Value  = 8
$ws.Cells.Item(828,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(828,3).Value  = "Coquimbo"
$ws.Cells.Item(828,4).Value  = 44826
$ws.Cells.Item(828,5).Value  = 4
$ws.Cells.Item(828,6).Value  = 100112004
$ws.Cells.Item(828,7).Value  = "Cebolla"
$ws.Cells.Item(828,8).Value  = "Sin especificar"
$ws.Cells.Item(828,9).Value  = "1a (guarda)"
$ws.Cells.Item(828,10).Value = 2600
$ws.Cells.Item(828,11).Value = 9300
$ws.Cells.Item(828,12).Value = 9500
$ws.Cells.Item(828,13).Value = 9400
$ws.Cells.Item(828,14).Value = "$/malla 16 kilos"
$ws.Cells.Item(828,15).Value = "Región de O'Higgins"
$ws.Cells.Item(828,16).Value = 588
$ws.Cells.Item(828,17).Value = 16
$ws.Cells.Item(828,18).Value = "Hortaliza"

# New row 829 data
$ws.Cells.Item(829,1).Value  = 8
$ws.Cells.Item(829,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(829,3).Value  = "Coquimbo"
$ws.Cells.Item(829,4).Value  = 44826
$ws.Cells.Item(829,5).Value  = 4
$ws.Cells.Item(829,6).Value  = 100112004
$ws.Cells.Item(829,7).Value  = "Cebolla"
$ws.Cells.Item(829,8).Value  = "Sin especificar"
$ws.Cells.Item(829,9).Value  = "2a (guarda)"
$ws.Cells.Item(829,10).Value = 1660
$ws.Cells.Item(829,11).Value = 9000
$ws.Cells.Item(829,12).Value = 9200
$ws.Cells.Item(829,13).Value = 9100
$ws.Cells.Item(829,14).Value = "$/malla 16 kilos"
$ws.Cells.Item(829,15).Value = "Región de O'Higgins"
$ws.Cells.Item(829,16).Value = 569
$ws.Cells.Item(829,17).Value = 16
$ws.Cells.Item(829,18).Value = "Hortaliza"
